# Add two new rows (51 and 52) of model-training results to Sheet1,
# extending the data range from A1:W50 to A1:W52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 51 ----
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 2).Value = 10
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0.0003
$ws.Cells.Item(51, 8).Value = "Regular"
$ws.Cells.Item(51, 12).Value = "<function relu at 0x11d707488>"
$ws.Cells.Item(51, 13).Value = 0.9351999759674072
$ws.Cells.Item(51, 14).Value = 0.3497999906539917
$ws.Cells.Item(51, 16).Value = 0.2422611862421036
$ws.Cells.Item(51, 17).Value = 3.15626335144043

# ---- Row 52 ----
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 30
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 0
$ws.Cells.Item(52, 5).Value = 0
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 0.0003
$ws.Cells.Item(52, 8).Value = "Regular"
$ws.Cells.Item(52, 12).Value = "<function relu at 0x11d707488>"
$ws.Cells.Item(52, 13).Value = 0.9556999802589417
$ws.Cells.Item(52, 14).Value = 0.04399999976158142
$ws.Cells.Item(52, 16).Value = 0.2216933816671371
$ws.Cells.Item(52, 17).Value = 23.85161018371582
$ws.Cells.Item(52, 20).Value = "weights/model_367.ckpt"
